# Sync attendance_reports: reverse the order of comma-separated "Recorded By"
# entries in column G (e.g. "System, user@example.com" -> "user@example.com, System").
# Cells with only a single entry (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val
    if ($text -notlike "*,*") {
        continue
    }

    $parts = $text -split ","
    $count = $parts.Count
    $reversed = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i].Trim()
    }

    $newVal = [string]::Join(", ", $reversed)
    $cell.Value = $newVal
}
